$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9044681191444397
$ws.Range("B1").Value = 1.564083099365234
$ws.Range("C1").Value = 4.351256370544434
$ws.Range("D1").Value = 2.271829128265381
$ws.Range("E1").Value = 1.483600854873657
